$wb = $excel.ActiveWorkbook

# Updated "want to go" counts (column F) for rows 2-9.
$updates = @{
    2 = 140
    3 = 1708
    4 = 30
    5 = 29
    6 = 474
    7 = 158
    8 = 78
    9 = 634
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
